# Updates the "cryptos" price/volume table with refreshed values.
# Note: for price cells (column D) whose new text looks like a plain
# number (e.g. "306.03"), a leading apostrophe is used to force Excel to
# store the value as text (preserving formatting such as trailing
# zeros), and the cell style is reset back to "Normal" right after so no
# stray quote-prefix / number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.261.78"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.269.47"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'306.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'97.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "'0.529"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "'6.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Value = "2.621.57"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "'14.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "2.284.54"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "42.130.56"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'12.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'68.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'238.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'23.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("D28").Value = "'37.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "'161.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'17.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0738"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("D37").Value = "'2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("D43").Value = "'19.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").Value = "1.942.43"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "'9.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").Value = "'53.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'92.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'71.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -2.04%  "
